# Generate Report for Handoff
# Refresh the "Latest Handoff Date(time)" column for every file row whose
# Status is "Ready for handoff" or "Handback transform failed" (i.e. rows
# that are being (re-)handed off as part of this report run). Rows that
# are already "Handed back" or still "In Translation" are left untouched.

$wb = $excel.ActiveWorkbook

# New handoff timestamps produced by this report run.
$tsOverview = "2016-03-24 00:30:54"
$tsZhCn     = "2016-03-24 00:30:49"
$tsDeDe     = "2016-03-24 00:30:54"

# --- Overview sheet: column D = "Latest Handoff Date" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$overviewRows = @(7, 10, 11, 12, 13, 14, 15, 16)
foreach ($r in $overviewRows) {
    $wsOverview.Cells.Item($r, 4).Value = $tsOverview
}

# --- zh-cn sheet: column E = "Latest Handoff Datetime" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$zhRows = @(7, 10, 11, 12, 13, 14, 15, 16)
foreach ($r in $zhRows) {
    $wsZhCn.Cells.Item($r, 5).Value = $tsZhCn
}

# --- de-de sheet: column E = "Latest Handoff Datetime" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$deRows = @(7, 10, 11, 12, 13, 14, 15, 16)
foreach ($r in $deRows) {
    $wsDeDe.Cells.Item($r, 5).Value = $tsDeDe
}
